$wb = $excel.ActiveWorkbook

# Sheet "About"
$about = $wb.Worksheets.Item("About")
$about.Range("A2").Value = "Version: mines - January 30 (built on February 02 2026 12.49.33 EST)"
$about.Range("A6").Value = "Recommended Citation:  `"Global Energy Monitor, Coal mine boundaries and methane sources for Xingwu Coal Mine, China, M2229, version 'mines - January 30 (built on February 02 2026 12.49.33 EST)'. (See the CC license for attribution requirements if sharing or adapting the data set.)"

# Sheet "Boundaries and methane sources"
$bms = $wb.Worksheets.Item("Boundaries and methane sources")
for ($r = 2; $r -le 10; $r++) {
    $bms.Cells.Item($r, 19).Value = "mines - January 30 (built on February 02 2026 12.49.33 EST)"
}
